$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (column headers) - drop unit suffixes / use short variable-style names
$ws.Range("A1").Value = "wojewodztwo"
$ws.Range("B1").Value = "lekarze"
$ws.Range("C1").Value = "pielegniarki"
$ws.Range("D1").Value = "apteki"
$ws.Range("E1").Value = "zgony_ogolem"
$ws.Range("F1").Value = "zespoly_ratownictwa"
$ws.Range("G1").Value = "ludnosc_na_lozko"
$ws.Range("H1").Value = "absolwenci"
$ws.Range("I1").Value = "organizacje_non_profit_ochrona_zdrowia"

# Update the selected cell in the sheet view to I1
$ws.Range("I1").Select()
